$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1602
$ws.Range("I32").Value = 1739
$ws.Range("J32").Value = 1562.8572
$ws.Range("K32").Value = 1739
$ws.Range("L32").Value = 1562.8572
$ws.Range("M32").Value = -1413
$ws.Range("N32").Value = -2214.8572
$ws.Range("H51").Value = 7300.3335
$ws.Range("J51").Value = 7300.3335
$ws.Range("L51").Value = 7300.3335
$ws.Range("N51").Value = -8268.3335
$ws.Range("H69").Value = 11589.454
$ws.Range("I69").Value = 8999.25
$ws.Range("K69").Value = 26997.75
$ws.Range("M69").Value = -26123.75
$ws.Range("H72").Value = 11589.454
$ws.Range("I72").Value = 8999.25
$ws.Range("K72").Value = 80993.25
$ws.Range("M72").Value = -76625.25
$ws.Range("H76").Value = 4138.4546
$ws.Range("I76").Value = 2841
$ws.Range("J76").Value = 5695.4
$ws.Range("K76").Value = 2841
$ws.Range("L76").Value = 5695.4
$ws.Range("M76").Value = -2526
$ws.Range("N76").Value = -6325.4
$ws.Range("H79").Value = 4138.4546
$ws.Range("I79").Value = 2841
$ws.Range("J79").Value = 5695.4
$ws.Range("K79").Value = 2841
$ws.Range("L79").Value = 5695.4
$ws.Range("M79").Value = -1749
$ws.Range("N79").Value = -7879.4
$ws.Range("H86").Value = 6745.706
$ws.Range("I86").Value = 1993.3334
$ws.Range("J86").Value = 7764.0713
$ws.Range("K86").Value = 1993.3334
$ws.Range("L86").Value = 7764.0713
$ws.Range("M86").Value = -870.3334
$ws.Range("N86").Value = -10010.0713
$ws.Range("H89").Value = 6745.706
$ws.Range("I89").Value = 1993.3334
$ws.Range("J89").Value = 7764.0713
$ws.Range("K89").Value = 9966.667
$ws.Range("L89").Value = 38820.35649999999
$ws.Range("M89").Value = -4350.666999999999
$ws.Range("N89").Value = -50052.35649999999
$ws.Range("H92").Value = 357.27777
$ws.Range("I92").Value = 350.53333
$ws.Range("K92").Value = 350.53333
$ws.Range("M92").Value = 897.46667
$ws.Range("H98").Value = 1688.3889
$ws.Range("I98").Value = 1092.3334
$ws.Range("K98").Value = 1092.3334
$ws.Range("M98").Value = 405.6666
$ws.Range("H106").Value = 27501926
$ws.Range("I106").Value = 33847624
$ws.Range("K106").Value = 33847624
$ws.Range("M106").Value = -33846993
$ws.Range("H113").Value = 5000
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -11508
$ws.Range("H116").Value = 4552.5713
$ws.Range("J116").Value = 4661.3335
$ws.Range("L116").Value = 4661.3335
$ws.Range("N116").Value = -11545.3335
$ws.Range("H122").Value = 1688.3889
$ws.Range("I122").Value = 1092.3334
$ws.Range("K122").Value = 3277.0002
$ws.Range("M122").Value = -827.0002
$ws.Range("H132").Value = 3098.5652
$ws.Range("I132").Value = 2858.825
$ws.Range("J132").Value = 4696.8335
$ws.Range("K132").Value = 8576.474999999999
$ws.Range("L132").Value = 14090.5005
$ws.Range("M132").Value = -6046.474999999999
$ws.Range("N132").Value = -19150.5005
$ws.Range("H137").Value = 6268.7617
$ws.Range("I137").Value = 4643.68
$ws.Range("J137").Value = 8658.588
$ws.Range("K137").Value = 13931.04
$ws.Range("L137").Value = 25975.764
$ws.Range("M137").Value = -11381.04
$ws.Range("N137").Value = -31075.764
$ws.Range("H138").Value = 4466.12
$ws.Range("I138").Value = 2587.5217
$ws.Range("J138").Value = 6066.407
$ws.Range("K138").Value = 7762.5651
$ws.Range("L138").Value = 18199.221
$ws.Range("M138").Value = -2622.5651
$ws.Range("N138").Value = -28479.221
$ws.Range("H141").Value = 5433.522
$ws.Range("I141").Value = 2269.6365
$ws.Range("J141").Value = 8333.75
$ws.Range("K141").Value = 6808.9095
$ws.Range("L141").Value = 25001.25
$ws.Range("M141").Value = -1628.9095
$ws.Range("N141").Value = -35361.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 129.125
$ws.Range("I4").Value = 122.57143
$ws.Range("K4").Value = 122.57143
$ws.Range("M4").Value = -6.571430000000007
$ws.Range("H32").Value = 4444.339
$ws.Range("I32").Value = 2914.6785
$ws.Range("K32").Value = 2914.6785
$ws.Range("M32").Value = -2627.6785
$ws.Range("H61").Value = 20838574
$ws.Range("I61").Value = 2521.9167
$ws.Range("J61").Value = 41674624
$ws.Range("K61").Value = 2521.9167
$ws.Range("L61").Value = 41674624
$ws.Range("M61").Value = -2309.9167
$ws.Range("N61").Value = -41675048
$ws.Range("H63").Value = 4385.5713
$ws.Range("I63").Value = 4000
$ws.Range("K63").Value = 4000
$ws.Range("M63").Value = -3314
$ws.Range("H66").Value = 4385.5713
$ws.Range("I66").Value = 4000
$ws.Range("K66").Value = 20000
$ws.Range("M66").Value = -16568
$ws.Range("H74").Value = 5429.2104
$ws.Range("I74").Value = 2836.4546
$ws.Range("K74").Value = 2836.4546
$ws.Range("M74").Value = -1962.4546
$ws.Range("H77").Value = 5429.2104
$ws.Range("I77").Value = 2836.4546
$ws.Range("K77").Value = 14182.273
$ws.Range("M77").Value = -9814.273000000001
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""
$ws.Range("H102").Value = 1605.5714
$ws.Range("I102").Value = 426.66666
$ws.Range("K102").Value = 426.66666
$ws.Range("M102").Value = 1195.33334
$ws.Range("H108").Value = 199950
$ws.Range("J108").Value = 199950
$ws.Range("L108").Value = 199950
$ws.Range("N108").Value = -207630
$ws.Range("H132").Value = 45457856
$ws.Range("I132").Value = 2085.8823
$ws.Range("J132").Value = 200007470
$ws.Range("K132").Value = 6257.646900000001
$ws.Range("L132").Value = 600022410
$ws.Range("M132").Value = -3727.646900000001
$ws.Range("N132").Value = -600027470
$ws.Range("H136").Value = 20838574
$ws.Range("I136").Value = 2521.9167
$ws.Range("J136").Value = 41674624
$ws.Range("K136").Value = 7565.750100000001
$ws.Range("L136").Value = 125023872
$ws.Range("M136").Value = -5015.750100000001
$ws.Range("N136").Value = -125028972

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = ""
$ws.Range("H21").Value = 54787.832
$ws.Range("J21").Value = 54787.832
$ws.Range("L21").Value = 54787.832
$ws.Range("N21").Value = -55259.832
$ws.Range("H34").Value = 10000
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10228
$ws.Range("H94").Value = 1408.3334
$ws.Range("I94").Value = 1264.5454
$ws.Range("J94").Value = 2990
$ws.Range("K94").Value = 1264.5454
$ws.Range("L94").Value = 2990
$ws.Range("M94").Value = -813.5454
$ws.Range("N94").Value = -3892
$ws.Range("H98").Value = 172593.5
$ws.Range("J98").Value = 172593.5
$ws.Range("L98").Value = 172593.5
$ws.Range("N98").Value = -178583.5
$ws.Range("H99").Value = 2223.3333
$ws.Range("I99").Value = 1499.8
$ws.Range("J99").Value = 3127.75
$ws.Range("K99").Value = 1499.8
$ws.Range("L99").Value = 3127.75
$ws.Range("M99").Value = -1.799999999999955
$ws.Range("N99").Value = -6123.75
$ws.Range("H103").Value = 145675.25
$ws.Range("J103").Value = 145675.25
$ws.Range("L103").Value = 145675.25
$ws.Range("N103").Value = -148019.25
$ws.Range("H105").Value = 2682.4211
$ws.Range("I105").Value = 2682.4211
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2682.4211
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -935.4211
$ws.Range("N105").Value = ""
$ws.Range("H107").Value = 71429290
$ws.Range("I107").Value = 459.8
$ws.Range("K107").Value = 459.8
$ws.Range("M107").Value = 1460.2
$ws.Range("H132").Value = 139999
$ws.Range("J132").Value = 139999
$ws.Range("L132").Value = 139999
$ws.Range("N132").Value = -150119
$ws.Range("H134").Value = 3202.3438
$ws.Range("I134").Value = 2215.8667
$ws.Range("K134").Value = 6647.6001
$ws.Range("M134").Value = -4112.6001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 252668.33
$ws.Range("H16").Value = 1477.5
$ws.Range("I16").Value = 1136.6666
$ws.Range("K16").Value = 1136.6666
$ws.Range("M16").Value = -849.6666
$ws.Range("H31").Value = 4738.203
$ws.Range("I31").Value = 2826.1428
$ws.Range("J31").Value = 9422.75
$ws.Range("K31").Value = 2826.1428
$ws.Range("L31").Value = 9422.75
$ws.Range("M31").Value = -2531.1428
$ws.Range("N31").Value = -10012.75
$ws.Range("H34").Value = 4738.203
$ws.Range("I34").Value = 2826.1428
$ws.Range("J34").Value = 9422.75
$ws.Range("K34").Value = 2826.1428
$ws.Range("L34").Value = 9422.75
$ws.Range("M34").Value = -2624.1428
$ws.Range("N34").Value = -9826.75
$ws.Range("H53").Value = 47966
$ws.Range("J53").Value = 47966
$ws.Range("L53").Value = 47966
$ws.Range("N53").Value = -49180
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = ""
$ws.Range("H86").Value = 6735.3335
$ws.Range("I86").Value = 5769.8335
$ws.Range("K86").Value = 5769.8335
$ws.Range("M86").Value = -4646.8335
$ws.Range("H89").Value = 6735.3335
$ws.Range("I89").Value = 5769.8335
$ws.Range("K89").Value = 28849.1675
$ws.Range("M89").Value = -23233.1675
$ws.Range("H96").Value = 126574.5
$ws.Range("J96").Value = 126574.5
$ws.Range("L96").Value = 126574.5
$ws.Range("N96").Value = -132066.5
$ws.Range("H105").Value = 52633004
$ws.Range("I105").Value = 1505.7222
$ws.Range("K105").Value = 1505.7222
$ws.Range("M105").Value = 241.2778000000001
$ws.Range("H107").Value = 588.2647
$ws.Range("I107").Value = 458.51852
$ws.Range("J107").Value = 1088.7142
$ws.Range("K107").Value = 458.51852
$ws.Range("L107").Value = 1088.7142
$ws.Range("M107").Value = 1461.48148
$ws.Range("N107").Value = -4928.7142
$ws.Range("H113").Value = 1477.5
$ws.Range("I113").Value = 1136.6666
$ws.Range("K113").Value = 1136.6666
$ws.Range("M113").Value = 1033.3334
$ws.Range("H122").Value = 1358.7142
$ws.Range("I122").Value = 1376
$ws.Range("J122").Value = 1255
$ws.Range("K122").Value = 4128
$ws.Range("L122").Value = 3765
$ws.Range("M122").Value = -1678
$ws.Range("N122").Value = -8665
$ws.Range("H132").Value = 5616
$ws.Range("J132").Value = 8504.667
$ws.Range("L132").Value = 25514.001
$ws.Range("N132").Value = -30574.001
$ws.Range("H134").Value = 7000.1816
$ws.Range("I134").Value = 6223
$ws.Range("J134").Value = 10497.5
$ws.Range("K134").Value = 18669
$ws.Range("L134").Value = 31492.5
$ws.Range("M134").Value = -16134
$ws.Range("N134").Value = -36562.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 136.28572
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 158.16667
$ws.Range("K7").Value = 15
$ws.Range("L7").Value = 474.50001
$ws.Range("M7").Value = 97
$ws.Range("N7").Value = -698.50001
$ws.Range("H12").Value = 357.94446
$ws.Range("I12").Value = 156.25
$ws.Range("J12").Value = 415.57144
$ws.Range("K12").Value = 468.75
$ws.Range("L12").Value = 1246.71432
$ws.Range("M12").Value = -295.75
$ws.Range("N12").Value = -1592.71432
$ws.Range("H14").Value = 1637.6666
$ws.Range("I14").Value = 1637.6666
$ws.Range("K14").Value = 4912.9998
$ws.Range("M14").Value = -4739.9998
$ws.Range("H23").Value = 45.666668
$ws.Range("J23").Value = 45.666668
$ws.Range("L23").Value = 137.000004
$ws.Range("N23").Value = -607.000004
$ws.Range("H32").Value = 1659.4
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1659.4
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 4978.200000000001
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -5544.200000000001
$ws.Range("H68").Value = 3653.348
$ws.Range("I68").Value = 2417.5557
$ws.Range("J68").Value = 4447.7856
$ws.Range("K68").Value = 7252.6671
$ws.Range("L68").Value = 13343.3568
$ws.Range("M68").Value = -6441.6671
$ws.Range("N68").Value = -14965.3568
$ws.Range("H71").Value = 3653.348
$ws.Range("I71").Value = 2417.5557
$ws.Range("J71").Value = 4447.7856
$ws.Range("K71").Value = 21758.0013
$ws.Range("L71").Value = 40030.0704
$ws.Range("M71").Value = -17702.0013
$ws.Range("N71").Value = -48142.0704
$ws.Range("H81").Value = 1992.6666
$ws.Range("I81").Value = 1992.6666
$ws.Range("K81").Value = 5977.9998
$ws.Range("M81").Value = -4854.9998
$ws.Range("H84").Value = 1992.6666
$ws.Range("I84").Value = 1992.6666
$ws.Range("K84").Value = 17933.9994
$ws.Range("M84").Value = -12317.9994
$ws.Range("H92").Value = 90
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H123").Value = 204959.8
$ws.Range("I123").Value = 204959.8
$ws.Range("K123").Value = 614879.3999999999
$ws.Range("M123").Value = -612429.3999999999
$ws.Range("H132").Value = 2622.4119
$ws.Range("I132").Value = 2435.875
$ws.Range("J132").Value = 2788.2222
$ws.Range("K132").Value = 21922.875
$ws.Range("L132").Value = 25093.9998
$ws.Range("M132").Value = -19392.875
$ws.Range("N132").Value = -30153.9998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 49999
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = ""
$ws.Range("H70").Value = 3802.158
$ws.Range("I70").Value = 2951.1538
$ws.Range("K70").Value = 2951.1538
$ws.Range("M70").Value = -2681.1538
$ws.Range("H73").Value = 3802.158
$ws.Range("I73").Value = 2951.1538
$ws.Range("K73").Value = 2951.1538
$ws.Range("M73").Value = -2015.1538
$ws.Range("H97").Value = 684.4286
$ws.Range("I97").Value = 464.8
$ws.Range("J97").Value = 1233.5
$ws.Range("K97").Value = 464.8
$ws.Range("L97").Value = 1233.5
$ws.Range("M97").Value = 31.19999999999999
$ws.Range("N97").Value = -2225.5
$ws.Range("H102").Value = 3072.5417
$ws.Range("I102").Value = 1648.0714
$ws.Range("J102").Value = 5066.8
$ws.Range("K102").Value = 1648.0714
$ws.Range("L102").Value = 5066.8
$ws.Range("M102").Value = -26.07140000000004
$ws.Range("N102").Value = -8310.8
$ws.Range("H113").Value = 2447.4119
$ws.Range("I113").Value = 1885.9048
$ws.Range("J113").Value = 3354.4614
$ws.Range("K113").Value = 1885.9048
$ws.Range("L113").Value = 3354.4614
$ws.Range("M113").Value = 284.0952
$ws.Range("N113").Value = -7694.4614
$ws.Range("H119").Value = 199941.67
$ws.Range("J119").Value = 199941.67
$ws.Range("L119").Value = 199941.67
$ws.Range("N119").Value = -209617.67
$ws.Range("H126").Value = 3535.5
$ws.Range("J126").Value = 5078
$ws.Range("L126").Value = 15234
$ws.Range("N126").Value = -20174
$ws.Range("H132").Value = 4028.7827
$ws.Range("I132").Value = 3441.5
$ws.Range("J132").Value = 7944
$ws.Range("K132").Value = 10324.5
$ws.Range("L132").Value = 23832
$ws.Range("M132").Value = -7794.5
$ws.Range("N132").Value = -28892

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1342.5555
$ws.Range("I22").Value = 1132.8334
$ws.Range("J22").Value = 1447.4166
$ws.Range("K22").Value = 1132.8334
$ws.Range("L22").Value = 1447.4166
$ws.Range("M22").Value = -837.8334
$ws.Range("N22").Value = -2037.4166
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = ""
$ws.Range("H27").Value = 1342.5555
$ws.Range("I27").Value = 1132.8334
$ws.Range("J27").Value = 1447.4166
$ws.Range("K27").Value = 1132.8334
$ws.Range("L27").Value = 1447.4166
$ws.Range("M27").Value = -1025.8334
$ws.Range("N27").Value = -1661.4166
$ws.Range("H40").Value = 6043.077
$ws.Range("I40").Value = 4869.091
$ws.Range("J40").Value = 12500
$ws.Range("K40").Value = 4869.091
$ws.Range("L40").Value = 12500
$ws.Range("M40").Value = -4733.091
$ws.Range("N40").Value = -12772
$ws.Range("H100").Value = 5551.273
$ws.Range("I100").Value = 5886.4
$ws.Range("J100").Value = 2200
$ws.Range("K100").Value = 5886.4
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -5345.4
$ws.Range("N100").Value = -3282
$ws.Range("H122").Value = 3940.238
$ws.Range("J122").Value = 5247.5
$ws.Range("L122").Value = 15742.5
$ws.Range("N122").Value = -20642.5
$ws.Range("H132").Value = 4755.4688
$ws.Range("I132").Value = 4635.875
$ws.Range("K132").Value = 13907.625
$ws.Range("M132").Value = -11377.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""
$ws.Range("H62").Value = 12400
$ws.Range("I62").Value = 12400
$ws.Range("K62").Value = 12400
$ws.Range("M62").Value = -11776
$ws.Range("H65").Value = 12400
$ws.Range("I65").Value = 12400
$ws.Range("K65").Value = 62000
$ws.Range("M65").Value = -58880
$ws.Range("H81").Value = 1430.75
$ws.Range("I81").Value = 1430.75
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2861.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1800.5
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 1430.75
$ws.Range("I84").Value = 1430.75
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 14307.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -9003.5
$ws.Range("N84").Value = ""
$ws.Range("H98").Value = 182750
$ws.Range("J98").Value = 182750
$ws.Range("L98").Value = 182750
$ws.Range("N98").Value = -188740
$ws.Range("H107").Value = 910598.8
$ws.Range("I107").Value = 1177855.8
$ws.Range("J107").Value = 1925.2
$ws.Range("K107").Value = 3533567.4
$ws.Range("L107").Value = 5775.6
$ws.Range("M107").Value = -3531647.4
$ws.Range("N107").Value = -9615.6
$ws.Range("H126").Value = 1940.6428
$ws.Range("I126").Value = 1816.9166
$ws.Range("J126").Value = 2683
$ws.Range("K126").Value = 5450.7498
$ws.Range("L126").Value = 8049
$ws.Range("M126").Value = -2980.7498
$ws.Range("N126").Value = -12989
$ws.Range("H132").Value = 3041.7812
$ws.Range("I132").Value = 2403
$ws.Range("K132").Value = 7209
$ws.Range("M132").Value = -4679
$ws.Range("H136").Value = 5304.984
$ws.Range("I136").Value = 4359.119
$ws.Range("J136").Value = 7196.7144
$ws.Range("K136").Value = 13077.357
$ws.Range("L136").Value = 21590.1432
$ws.Range("M136").Value = -10527.357
$ws.Range("N136").Value = -26690.1432
